{"js": "// Ghost Slider review: drop the standalone \"Meta description\" paragraph\n// near the top, and move its description text down to the bottom of the\n// document (replacing the old AI image-generation \"Prompt:\" paragraph),\n// preceded by a new bold paragraph repeating the page title.\n\nconst body = context.document.body;\n\n// --- Step 1: remove the \"Meta description\" paragraph (2nd paragraph). ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst metaDescriptionParagraph = paragraphs.items[1];\nmetaDescriptionParagraph.delete();\nawait context.sync();\n\n// --- Step 2: insert a new bold paragraph with the page title, right    ---\n// --- before the final \"Prompt: ...\" paragraph.                        ---\nconst remainingParagraphs = body.paragraphs;\nremainingParagraphs.load(\"text\");\nawait context.sync();\n\nconst promptParagraph = remainingParagraphs.items[remainingParagraphs.items.length - 1];\n\nconst titleParagraph = promptParagraph.insertParagraph(\"\", \"Before\");\n// Start from a clean (unformatted) paragraph so the new text doesn't\n// inherit the italic formatting of the neighboring \"Prompt:\" paragraph.\ntitleParagraph.clear();\nconst titleRange = titleParagraph.getRange();\ntitleRange.insertText(\n  \"Play Ghost Slider for Free: Unique Free Spin Mode for Big Payouts\",\n  \"Replace\"\n);\ntitleRange.font.bold = true;\nawait context.sync();\n\n// --- Step 3: replace the \"Prompt: ...\" paragraph's text with the new   ---\n// --- meta-description text; the matched run's existing italic          ---\n// --- formatting is preserved automatically.                            ---\nconst oldPromptText =\n  \"Prompt: Create a featured image for Ghost Slider that captures the spooky adventure and features a happy Maya warrior with glasses. The image should be in a cartoon style to match the game's design and feature vibrant colors. The warrior should have a confident and adventurous look on their face, maybe holding a torch or pointing a sword towards the ghosts around them. The background should feature a graveyard or haunted mansion with ghostly apparitions lurking in the shadows. The overall tone should be spooky but with a touch of fun and adventure, inviting players to join the warrior in their ghostly quest.\";\nconst newDescriptionText =\n  \"Read our review of Ghost Slider, a 5*3 slot with a ghost theme and a unique free spin mode for experienced players seeking big payouts. Play for free now.\";\n\nconst searchResults = body.search(oldPromptText, { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nsearchResults.items[0].insertText(newDescriptionText, \"Replace\");\nawait context.sync();\n", "ps1": "# Ghost Slider review: drop the standalone \"Meta description\" paragraph\n# near the top, and move its description text down to the bottom of the\n# document (replacing the old AI image-generation \"Prompt:\" paragraph),\n# preceded by a new bold paragraph repeating the page title.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: remove the \"Meta description\" paragraph (2nd paragraph). ---\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Range.Delete()\n\n# --- Step 2: replace the \"Prompt: ...\" paragraph's text with the new   ---\n# --- meta-description text. Find/Replace keeps the run's existing      ---\n# --- italic formatting intact.                                         ---\n$promptText = \"Prompt: Create a featured image for Ghost Slider that captures the spooky adventure and features a happy Maya warrior with glasses. The image should be in a cartoon style to match the game's design and feature vibrant colors. The warrior should have a confident and adventurous look on their face, maybe holding a torch or pointing a sword towards the ghosts around them. The background should feature a graveyard or haunted mansion with ghostly apparitions lurking in the shadows. The overall tone should be spooky but with a touch of fun and adventure, inviting players to join the warrior in their ghostly quest.\"\n$newDescText = \"Read our review of Ghost Slider, a 5*3 slot with a ghost theme and a unique free spin mode for experienced players seeking big payouts. Play for free now.\"\n\n$findRange = $d.Content\n$findRange.Find.Execute($promptText, $false, $false, $false, $false, $false, $true, 1, $false, $newDescText, 2) | Out-Null\n\n# --- Step 3: insert a new paragraph just before that one, with the     ---\n# --- page title text in bold, matching the \"Normal\" style used by the  ---\n# --- rest of the body text.                                            ---\n$anchorIndex = $d.Paragraphs.Count - 1\n$anchorPara = $d.Paragraphs.Item($anchorIndex)\n$anchorPara.Range.InsertParagraphAfter()\n\n$titleParaIndex = $anchorIndex + 1\n$titlePara = $d.Paragraphs.Item($titleParaIndex)\n$titlePara.Style = \"Normal\"\n$titlePara.Range.Text = \"Play Ghost Slider for Free: Unique Free Spin Mode for Big Payouts\"\n\n$titlePara = $d.Paragraphs.Item($titleParaIndex)\n$titleRange = $titlePara.Range\n$titleRange.MoveEnd(1, -1) | Out-Null\n$titleRange.Font.Bold = 1\n"}
